$d = $word.ActiveDocument

# Locate the end of the title line "... Malayalam Corrections - prior to 31st Aug 2022"
$rng = $d.Content.Duplicate
$rng.Find.Execute("prior to 31st Aug 2022", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)

# Select the collapsed insertion point and type the missing period so a brand
# new run is created (matching the author's edit) instead of being merged
# into the preceding run.
$rng.Select()
$sel = $word.Selection
$sel.TypeText(".")

# Give the newly typed run the same character formatting as the rest of the
# heading (bold, bold-complex-script, 16pt/16pt-cs, single underline).
$sel.Font.Bold = $true
$sel.Font.BoldBi = $true
$sel.Font.Size = 16
$sel.Font.SizeBi = 16
$sel.Font.Underline = 1
